$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$question = "How many tracks can you define in one ODF?"
$answer = "According to the Track Settings dialog box, the number of tracks that can be defined is 200."

$rows = @(
    @{ Row = 7; Model = "llama3.2:latest" },
    @{ Row = 8; Model = "deepseek1.5" },
    @{ Row = 9; Model = "openai" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $question
    $ws.Cells.Item($r.Row, 2).Value = $r.Model
    $ws.Cells.Item($r.Row, 3).Value = $answer
}
